$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell "intervention_type" in column K, matching the style used
# by the rest of the header row (J1 -> K1)
$ws.Cells.Item(1, 11).Value = "intervention_type"
$ws.Cells.Item(1, 10).Copy()
$ws.Cells.Item(1, 11).PasteSpecial(-4122)

# Values for the new "intervention_type" column, rows 2-34
$values = @(
    "DRUG",
    "PROCEDURE",
    "PROCEDURE",
    "OTHER",
    "DRUG",
    "OTHER",
    "OTHER",
    "DEVICE",
    "OTHER",
    "DEVICE",
    "DEVICE",
    "DEVICE",
    "OTHER",
    "BIOLOGICAL",
    "DEVICE",
    "PROCEDURE",
    "BIOLOGICAL",
    "OTHER",
    "DEVICE",
    "DRUG",
    "OTHER",
    "PROCEDURE",
    "OTHER",
    "DRUG",
    "DEVICE",
    "OTHER",
    "PROCEDURE",
    "DIAGNOSTIC_TEST",
    "PROCEDURE",
    "DEVICE",
    "OTHER",
    "DEVICE",
    "DRUG"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $values[$i]
}
